# Rename the table header cells from Finnish/Japanese ("Japani"/"Suomi")
# to English ("japanese"/"finnish"), and move the active selection to B2
# (matching the author's commit: the first-row headers were translated
# and the sheet's saved cursor position moved off the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "japanese"
$ws.Range("B1").Value = "finnish"

$ws.Range("B2").Select()
